$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update F9: "Yes" -> "No" (do first so "No" becomes shared-string index 27) ---
$ws.Range("F9").Value = "No"

# --- Update D3: "Nil" -> URL, with new bold blue Courier New font + vertical-center alignment ---
$d3 = $ws.Range("D3")
$d3.Value = "https://opensource-demo.orangehrmlive.com/"
$d3.Font.Bold = $true
$d3.Font.Size = 10
$d3.Font.Color = 16711722
$d3.Font.Name = "Courier New"
$d3.Font.Family = 3
$d3.VerticalAlignment = -4108

# --- Add new row 10 (all values already exist as shared strings) ---
$ws.Range("C10").Value = "Click"
$ws.Range("D10").Value = "Nil"
$ws.Range("E10").Value = "Admin"
$ws.Range("F10").Value = "Yes"

# --- Add new row 11 ---
$ws.Range("C11").Value = "Click"
$ws.Range("D11").Value = "Nil"
$ws.Range("E11").Value = "Add_Button"
$ws.Range("F11").Value = "Yes"

# --- Add new row 12 (fill E before C to match original shared-string insertion order) ---
$ws.Range("D12").Value = "Admin"
$ws.Range("E12").Value = "UserRoleSelection"
$ws.Range("C12").Value = "Dropdown"
$ws.Range("F12").Value = "Yes"

# --- Resize columns D and E (closest achievable values given column width quantization) ---
$ws.Columns("D").ColumnWidth = 50.709635416666664
$ws.Columns("E").ColumnWidth = 16.709635416666668

# --- Update selected cell to B12 ---
$ws.Range("B12").Select() | Out-Null
